# edit.ps1 - apply the two textual edits described by the diff:
#  1) "  Carl August of Saxe-Weimar grants the first German constitution."
#     becomes five runs:
#        "  " / "Grand Duke " / "Carl August of Saxe-Weimar" / "-Eisenach" / " grants the first German constitution."
#  2) "January 2016" becomes two runs: "May" / " 2016"

$d = $word.ActiveDocument

$wpNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-RangeWithRuns {
    param([string]$SearchText, [string]$RunsXml)

    foreach ($p in $d.Paragraphs) {
        $paraText = $p.Range.Text
        if ($paraText.Contains($SearchText)) {
            $seek = $p.Range.Duplicate
            $found = $seek.Find.Execute($SearchText, $true, $false, $false, $false, $false, `
                                         $true, 1, $false, "", 0)
            if ($found) {
                $target = $d.Range($seek.Start, $seek.End)
                $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
                       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
                       '<pkg:xmlData><w:document ' + $wpNs + '><w:body><w:p>' + $RunsXml + '</w:p></w:body></w:document></pkg:xmlData>' + `
                       '</pkg:part></pkg:package>'
                $target.InsertXML($xml)
                break
            }
        }
    }
}

# --- Edit 1 ---------------------------------------------------------------
$runs1 = '<w:r w:rsidRPr="00813531"><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r>' + `
         '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Grand Duke </w:t></w:r>' + `
         '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Carl August of Saxe-Weimar</w:t></w:r>' + `
         '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>-Eisenach</w:t></w:r>' + `
         '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> grants the first German constitution.</w:t></w:r>'

Replace-RangeWithRuns "  Carl August of Saxe-Weimar grants the first German constitution." $runs1

# --- Edit 2 -----------------------------------------------------------------
$runs2 = '<w:r><w:t>May</w:t></w:r>' + `
         '<w:r><w:t xml:space="preserve"> 2016</w:t></w:r>'

Replace-RangeWithRuns "January 2016" $runs2

Write-Host "Done."
